$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Misc")

$rows = @(
    @("M-12V-1A-Lead-Battery-Charger", "Part", "https://nl.rs-online.com/web/p/battery-chargers/4784077"),
    @("M-12V-Connector", "Part", "size reference: https://nl.rs-online.com/web/p/dc-power-connectors/8787203"),
    @("M-Emergency-Switch", "Part", "https://www.conrad.nl/nl/p/tru-components-las1-by-11tsb-noodstopschakelaar-250-v-ac-3-a-1x-nc-1x-no-ip40-1-stuk-s-704766.html?refresh=true#productTechData"),
    @("M-On-Off-Switch", "Part", "https://www.conrad.nl/nl/p/arcolectric-bulgin-ltd-c6053alnah-wipschakelaar-c6053alnah-230-v-ac-16-a-2x-uit-aan-continu-1-stuk-s-701528.html")
)

$startRow = 15
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}

$ws.Columns.Item(1).ColumnWidth = 29.166666666666668
$ws.Range("C20").Select()
